$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 47.0283035
$ws.Range("H2").Value = 94.056607
$ws.Range("I2").Value = 0.1170896029811303
$ws.Range("J2").Value = 0.08213676148506427
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.4088795
$ws.Range("N2").Value = 6.817759000000001
$ws.Range("O2").Value = 0.04462917857418398
$ws.Range("P2").Value = 0.03188023676415588
$ws.Range("Q2").Value = 160.3138197209283
$ws.Range("R2").Value = 641.255278883713
$ws.Range("S2").Value = 0.00522561280062517
$ws.Range("T2").Value = 0.002618539403184848

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 47.0283035
$ws.Range("H3").Value = 94.056607
$ws.Range("I3").Value = 0.1170896029811303
$ws.Range("J3").Value = 0.08213676148506427
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.142411
$ws.Range("N3").Value = 129.427233
$ws.Range("O3").Value = 0.5648220667934549
$ws.Range("P3").Value = 0.6052092530360149
$ws.Range("Q3").Value = 2028.914398229739
$ws.Range("R3").Value = 12173.48638937843
$ws.Range("S3").Value = 0.06613479155582712
$ws.Range("T3").Value = 0.04970992806517307

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 47.0283035
$ws.Range("H4").Value = 94.056607
$ws.Range("I4").Value = 0.1170896029811303
$ws.Range("J4").Value = 0.08213676148506427
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.580051
$ws.Range("N4").Value = 16.740153
$ws.Range("O4").Value = 0.07305423747951603
$ws.Range("P4").Value = 0.07827792697104638
$ws.Range("Q4").Value = 262.4203319734785
$ws.Range("R4").Value = 1574.521991840871
$ws.Range("S4").Value = 0.008553891662565743
$ws.Range("T4").Value = 0.006429495417166116

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 47.0283035
$ws.Range("H5").Value = 94.056607
$ws.Range("I5").Value = 0.1170896029811303
$ws.Range("J5").Value = 0.08213676148506427
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.215490333333333
$ws.Range("N5").Value = 24.646471
$ws.Range("O5").Value = 0.1075575083134548
$ws.Range("P5").Value = 0.1152483287955619
$ws.Range("Q5").Value = 386.3605727973161
$ws.Range("R5").Value = 2318.163436783897
$ws.Range("S5").Value = 0.01259386594606204
$ws.Range("T5").Value = 0.009466124493833333

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 47.0283035
$ws.Range("H6").Value = 94.056607
$ws.Range("I6").Value = 0.1170896029811303
$ws.Range("J6").Value = 0.08213676148506427
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.152783
$ws.Range("N6").Value = 12.458349
$ws.Range("O6").Value = 0.05436839116396912
$ws.Range("P6").Value = 0.05825596296532109
$ws.Range("Q6").Value = 195.2983392936405
$ws.Range("R6").Value = 1171.790035761843
$ws.Range("S6").Value = 0.006365973336111938
$ws.Range("T6").Value = 0.004784956135165316

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 47.0283035
$ws.Range("H7").Value = 94.056607
$ws.Range("I7").Value = 0.1170896029811303
$ws.Range("J7").Value = 0.08213676148506427
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.8826895
$ws.Range("N7").Value = 23.765379
$ws.Range("O7").Value = 0.1555686176754212
$ws.Range("P7").Value = 0.1111282914678999
$ws.Range("Q7").Value = 558.8227282022632
$ws.Range("R7").Value = 2235.290912809053
$ws.Range("S7").Value = 0.01821546767993832
$ws.Range("T7").Value = 0.009127717970541599

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.33126566666667
$ws.Range("H8").Value = 36.993797
$ws.Range("I8").Value = 0.03070200057641551
$ws.Range("J8").Value = 0.03230555276798244
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.4088795
$ws.Range("N8").Value = 6.817759000000001
$ws.Range("O8").Value = 0.04462917857418398
$ws.Range("P8").Value = 0.03188023676415588
$ws.Range("Q8").Value = 42.03579874015384
$ws.Range("R8").Value = 252.214792440923
$ws.Range("S8").Value = 0.001370205066309547
$ws.Range("T8").Value = 0.001029908671040211

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.33126566666667
$ws.Range("H9").Value = 36.993797
$ws.Range("I9").Value = 0.03070200057641551
$ws.Range("J9").Value = 0.03230555276798244
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 43.142411
$ws.Range("N9").Value = 129.427233
$ws.Range("O9").Value = 0.5648220667934549
$ws.Range("P9").Value = 0.6052092530360149
$ws.Range("Q9").Value = 532.0005315415224
$ws.Range("R9").Value = 4788.004783873701
$ws.Range("S9").Value = 0.01734116742026485
$ws.Range("T9").Value = 0.01955161945962621

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.33126566666667
$ws.Range("H10").Value = 36.993797
$ws.Range("I10").Value = 0.03070200057641551
$ws.Range("J10").Value = 0.03230555276798244
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.580051
$ws.Range("N10").Value = 16.740153
$ws.Range("O10").Value = 0.07305423747951603
$ws.Range("P10").Value = 0.07827792697104638
$ws.Range("Q10").Value = 68.80909131454901
$ws.Range("R10").Value = 619.281821830941
$ws.Range("S10").Value = 0.002242911241205696
$ws.Range("T10").Value = 0.002528811700331415

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 12.33126566666667
$ws.Range("H11").Value = 36.993797
$ws.Range("I11").Value = 0.03070200057641551
$ws.Range("J11").Value = 0.03230555276798244
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.215490333333333
$ws.Range("N11").Value = 24.646471
$ws.Range("O11").Value = 0.1075575083134548
$ws.Range("P11").Value = 0.1152483287955619
$ws.Range("Q11").Value = 101.3073938822652
$ws.Range("R11").Value = 911.766544940387
$ws.Range("S11").Value = 0.003302230682237504
$ws.Range("T11").Value = 0.003723160967326815

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 12.33126566666667
$ws.Range("H12").Value = 36.993797
$ws.Range("I12").Value = 0.03070200057641551
$ws.Range("J12").Value = 0.03230555276798244
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.152783
$ws.Range("N12").Value = 12.458349
$ws.Range("O12").Value = 0.05436839116396912
$ws.Range("P12").Value = 0.05825596296532109
$ws.Range("Q12").Value = 51.209070429017
$ws.Range("R12").Value = 460.881633861153
$ws.Range("S12").Value = 0.001669218376854964
$ws.Range("T12").Value = 0.001881991085625811

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 12.33126566666667
$ws.Range("H13").Value = 36.993797
$ws.Range("I13").Value = 0.03070200057641551
$ws.Range("J13").Value = 0.03230555276798244
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 11.8826895
$ws.Range("N13").Value = 23.765379
$ws.Range("O13").Value = 0.1555686176754212
$ws.Range("P13").Value = 0.1111282914678999
$ws.Range("Q13").Value = 146.5286010590105
$ws.Range("R13").Value = 879.171606354063
$ws.Range("S13").Value = 0.004776267789542945
$ws.Range("T13").Value = 0.003590060884031973

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 129.0494106666667
$ws.Range("H14").Value = 387.148232
$ws.Range("I14").Value = 0.3213031969122349
$ws.Range("J14").Value = 0.3380847237148192
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.4088795
$ws.Range("N14").Value = 6.817759000000001
$ws.Range("O14").Value = 0.04462917857418398
$ws.Range("P14").Value = 0.03188023676415588
$ws.Range("Q14").Value = 439.9138905086813
$ws.Range("R14").Value = 2639.483343052088
$ws.Range("S14").Value = 0.01433949775145233
$ws.Range("T14").Value = 0.01077822103837266

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 129.0494106666667
$ws.Range("H15").Value = 387.148232
$ws.Range("I15").Value = 0.3213031969122349
$ws.Range("J15").Value = 0.3380847237148192
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 43.142411
$ws.Range("N15").Value = 129.427233
$ws.Range("O15").Value = 0.5648220667934549
$ws.Range("P15").Value = 0.6052092530360149
$ws.Range("Q15").Value = 5567.502714289118
$ws.Range("R15").Value = 50107.52442860206
$ws.Range("S15").Value = 0.1814791357473129
$ws.Range("T15").Value = 0.2046120031023332

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 129.0494106666667
$ws.Range("H16").Value = 387.148232
$ws.Range("I16").Value = 0.3213031969122349
$ws.Range("J16").Value = 0.3380847237148192
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 5.580051
$ws.Range("N16").Value = 16.740153
$ws.Range("O16").Value = 0.07305423747951603
$ws.Range("P16").Value = 0.07827792697104638
$ws.Range("Q16").Value = 720.102293039944
$ws.Range("R16").Value = 6480.920637359496
$ws.Range("S16").Value = 0.02347256005015411
$ws.Range("T16").Value = 0.02646457131297501

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 129.0494106666667
$ws.Range("H17").Value = 387.148232
$ws.Range("I17").Value = 0.3213031969122349
$ws.Range("J17").Value = 0.3380847237148192
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 8.215490333333333
$ws.Range("N17").Value = 24.646471
$ws.Range("O17").Value = 0.1075575083134548
$ws.Range("P17").Value = 0.1152483287955619
$ws.Range("Q17").Value = 1060.204185854363
$ws.Range("R17").Value = 9541.837672689271
$ws.Range("S17").Value = 0.0345585712730273
$ws.Range("T17").Value = 0.03896369939944219

$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 129.0494106666667
$ws.Range("H18").Value = 387.148232
$ws.Range("I18").Value = 0.3213031969122349
$ws.Range("J18").Value = 0.3380847237148192
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 4.152783
$ws.Range("N18").Value = 12.458349
$ws.Range("O18").Value = 0.05436839116396912
$ws.Range("P18").Value = 0.05825596296532109
$ws.Range("Q18").Value = 535.914198776552
$ws.Range("R18").Value = 4823.227788988968
$ws.Range("S18").Value = 0.01746873789195818
$ws.Range("T18").Value = 0.01969545114387132

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 129.0494106666667
$ws.Range("H19").Value = 387.148232
$ws.Range("I19").Value = 0.3213031969122349
$ws.Range("J19").Value = 0.3380847237148192
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 11.8826895
$ws.Range("N19").Value = 23.765379
$ws.Range("O19").Value = 0.1555686176754212
$ws.Range("P19").Value = 0.1111282914678999
$ws.Range("Q19").Value = 1533.454077109988
$ws.Range("R19").Value = 9200.724462659928
$ws.Range("S19").Value = 0.04998469419833003
$ws.Range("T19").Value = 0.03757077771782484

$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 158.6435343333333
$ws.Range("H20").Value = 475.930603
$ws.Range("I20").Value = 0.3949857228129294
$ws.Range("J20").Value = 0.4156156560277983
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 3.4088795
$ws.Range("N20").Value = 6.817759000000001
$ws.Range("O20").Value = 0.04462917857418398
$ws.Range("P20").Value = 0.03188023676415588
$ws.Range("Q20").Value = 540.7966919964462
$ws.Range("R20").Value = 3244.780151978677
$ws.Range("S20").Value = 0.01762788835767136
$ws.Range("T20").Value = 0.01324992551705618

$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 158.6435343333333
$ws.Range("H21").Value = 475.930603
$ws.Range("I21").Value = 0.3949857228129294
$ws.Range("J21").Value = 0.4156156560277983
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 43.142411
$ws.Range("N21").Value = 129.427233
$ws.Range("O21").Value = 0.5648220667934549
$ws.Range("P21").Value = 0.6052092530360149
$ws.Range("Q21").Value = 6844.264560701277
$ws.Range("R21").Value = 61598.38104631149
$ws.Range("S21").Value = 0.2230966523131055
$ws.Range("T21").Value = 0.2515344407346571

$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 158.6435343333333
$ws.Range("H22").Value = 475.930603
$ws.Range("I22").Value = 0.3949857228129294
$ws.Range("J22").Value = 0.4156156560277983
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 5.580051
$ws.Range("N22").Value = 16.740153
$ws.Range("O22").Value = 0.07305423747951603
$ws.Range("P22").Value = 0.07827792697104638
$ws.Range("Q22").Value = 885.239012400251
$ws.Range("R22").Value = 7967.151111602258
$ws.Range("S22").Value = 0.02885538079539404
$ws.Range("T22").Value = 0.03253353197056753

$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 158.6435343333333
$ws.Range("H23").Value = 475.930603
$ws.Range("I23").Value = 0.3949857228129294
$ws.Range("J23").Value = 0.4156156560277983
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 8.215490333333333
$ws.Range("N23").Value = 24.646471
$ws.Range("O23").Value = 0.1075575083134548
$ws.Range("P23").Value = 0.1152483287955619
$ws.Range("Q23").Value = 1303.334422761335
$ws.Range("R23").Value = 11730.00980485201
$ws.Range("S23").Value = 0.04248368016514759
$ws.Range("T23").Value = 0.04789900977847487

$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 158.6435343333333
$ws.Range("H24").Value = 475.930603
$ws.Range("I24").Value = 0.3949857228129294
$ws.Range("J24").Value = 0.4156156560277983
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 4.152783
$ws.Range("N24").Value = 12.458349
$ws.Range("O24").Value = 0.05436839116396912
$ws.Range("P24").Value = 0.05825596296532109
$ws.Range("Q24").Value = 658.812172439383
$ws.Range("R24").Value = 5929.309551954447
$ws.Range("S24").Value = 0.02147473828207642
$ws.Range("T24").Value = 0.02421209026536305

$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 158.6435343333333
$ws.Range("H25").Value = 475.930603
$ws.Range("I25").Value = 0.3949857228129294
$ws.Range("J25").Value = 0.4156156560277983
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 11.8826895
$ws.Range("N25").Value = 23.765379
$ws.Range("O25").Value = 0.1555686176754212
$ws.Range("P25").Value = 0.1111282914678999
$ws.Range("Q25").Value = 1885.111859665589
$ws.Range("R25").Value = 11310.67115799354
$ws.Range("S25").Value = 0.0614473828995345
$ws.Range("T25").Value = 0.04618665776167961

$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 41.810285
$ws.Range("H26").Value = 125.430855
$ws.Range("I26").Value = 0.1040979432987182
$ws.Range("J26").Value = 0.1095349337872956
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 3.4088795
$ws.Range("N26").Value = 6.817759000000001
$ws.Range("O26").Value = 0.04462917857418398
$ws.Range("P26").Value = 0.03188023676415588
$ws.Range("Q26").Value = 142.5262234256575
$ws.Range("R26").Value = 855.1573405539451
$ws.Range("S26").Value = 0.004645805700683771
$ws.Range("T26").Value = 0.00349199962308512

$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 41.810285
$ws.Range("H27").Value = 125.430855
$ws.Range("I27").Value = 0.1040979432987182
$ws.Range("J27").Value = 0.1095349337872956
$ws.Range("K27").Value = 3
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = 43.142411
$ws.Range("N27").Value = 129.427233
$ws.Range("O27").Value = 0.5648220667934549
$ws.Range("P27").Value = 0.6052092530360149
$ws.Range("Q27").Value = 1803.796499497135
$ws.Range("R27").Value = 16234.16849547422
$ws.Range("S27").Value = 0.05879681548292987
$ws.Range("T27").Value = 0.0662915554587585

$ws.Range("E28").Value = 3
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 41.810285
$ws.Range("H28").Value = 125.430855
$ws.Range("I28").Value = 0.1040979432987182
$ws.Range("J28").Value = 0.1095349337872956
$ws.Range("K28").Value = 3
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 5.580051
$ws.Range("N28").Value = 16.740153
$ws.Range("O28").Value = 0.07305423747951603
$ws.Range("P28").Value = 0.07827792697104638
$ws.Range("Q28").Value = 233.303522624535
$ws.Range("R28").Value = 2099.731703620815
$ws.Range("S28").Value = 0.007604795870873751
$ws.Range("T28").Value = 0.008574167547780323

$ws.Range("E29").Value = 3
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 41.810285
$ws.Range("H29").Value = 125.430855
$ws.Range("I29").Value = 0.1040979432987182
$ws.Range("J29").Value = 0.1095349337872956
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = 1
$ws.Range("M29").Value = 8.215490333333333
$ws.Range("N29").Value = 24.646471
$ws.Range("O29").Value = 0.1075575083134548
$ws.Range("P29").Value = 0.1152483287955619
$ws.Range("Q29").Value = 343.4919922514117
$ws.Range("R29").Value = 3091.427930262705
$ws.Range("S29").Value = 0.01119651540176542
$ws.Range("T29").Value = 0.01262371806371835

$ws.Range("E30").Value = 3
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 41.810285
$ws.Range("H30").Value = 125.430855
$ws.Range("I30").Value = 0.1040979432987182
$ws.Range("J30").Value = 0.1095349337872956
$ws.Range("K30").Value = 3
$ws.Range("L30").Value = 1
$ws.Range("M30").Value = 4.152783
$ws.Range("N30").Value = 12.458349
$ws.Range("O30").Value = 0.05436839116396912
$ws.Range("P30").Value = 0.05825596296532109
$ws.Range("Q30").Value = 173.629040773155
$ws.Range("R30").Value = 1562.661366958395
$ws.Range("S30").Value = 0.005659637700629386
$ws.Range("T30").Value = 0.006381063046121589

$ws.Range("E31").Value = 3
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 41.810285
$ws.Range("H31").Value = 125.430855
$ws.Range("I31").Value = 0.1040979432987182
$ws.Range("J31").Value = 0.1095349337872956
$ws.Range("K31").Value = 2
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 11.8826895
$ws.Range("N31").Value = 23.765379
$ws.Range("O31").Value = 0.1555686176754212
$ws.Range("P31").Value = 0.1111282914678999
$ws.Range("Q31").Value = 496.8186345615075
$ws.Range("R31").Value = 2980.911807369045
$ws.Range("S31").Value = 0.01619437314183596
$ws.Range("T31").Value = 0.0121724300478317

$ws.Range("E32").Value = 2
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 12.780919
$ws.Range("H32").Value = 25.561838
$ws.Range("I32").Value = 0.03182153341857176
$ws.Range("J32").Value = 0.02232237221704002
$ws.Range("K32").Value = 2
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = 3.4088795
$ws.Range("N32").Value = 6.817759000000001
$ws.Range("O32").Value = 0.04462917857418398
$ws.Range("P32").Value = 0.03188023676415588
$ws.Range("Q32").Value = 43.56861277026051
$ws.Range("R32").Value = 174.274451081042
$ws.Range("S32").Value = 0.001420168897441802
$ws.Range("T32").Value = 0.0007116425114168511

$ws.Range("E33").Value = 2
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 12.780919
$ws.Range("H33").Value = 25.561838
$ws.Range("I33").Value = 0.03182153341857176
$ws.Range("J33").Value = 0.02232237221704002
$ws.Range("K33").Value = 3
$ws.Range("L33").Value = 1
$ws.Range("M33").Value = 43.142411
$ws.Range("N33").Value = 129.427233
$ws.Range("O33").Value = 0.5648220667934549
$ws.Range("P33").Value = 0.6052092530360149
$ws.Range("Q33").Value = 551.3996604557091
$ws.Range("R33").Value = 3308.397962734254
$ws.Range("S33").Value = 0.0179735042740147
$ws.Range("T33").Value = 0.01350970621546669

$ws.Range("E34").Value = 2
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 12.780919
$ws.Range("H34").Value = 25.561838
$ws.Range("I34").Value = 0.03182153341857176
$ws.Range("J34").Value = 0.02232237221704002
$ws.Range("K34").Value = 3
$ws.Range("L34").Value = 1
$ws.Range("M34").Value = 5.580051
$ws.Range("N34").Value = 16.740153
$ws.Range("O34").Value = 0.07305423747951603
$ws.Range("P34").Value = 0.07827792697104638
$ws.Range("Q34").Value = 71.318179846869
$ws.Range("R34").Value = 427.909079081214
$ws.Range("S34").Value = 0.002324697859322697
$ws.Range("T34").Value = 0.001747349022225974

$ws.Range("E35").Value = 2
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 12.780919
$ws.Range("H35").Value = 25.561838
$ws.Range("I35").Value = 0.03182153341857176
$ws.Range("J35").Value = 0.02232237221704002
$ws.Range("K35").Value = 3
$ws.Range("L35").Value = 1
$ws.Range("M35").Value = 8.215490333333333
$ws.Range("N35").Value = 24.646471
$ws.Range("O35").Value = 0.1075575083134548
$ws.Range("P35").Value = 0.1152483287955619
$ws.Range("Q35").Value = 105.0015164956163
$ws.Range("R35").Value = 630.009098973698
$ws.Range("S35").Value = 0.003422644845214911
$ws.Range("T35").Value = 0.002572616092766345

$ws.Range("E36").Value = 2
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 12.780919
$ws.Range("H36").Value = 25.561838
$ws.Range("I36").Value = 0.03182153341857176
$ws.Range("J36").Value = 0.02232237221704002
$ws.Range("K36").Value = 3
$ws.Range("L36").Value = 1
$ws.Range("M36").Value = 4.152783
$ws.Range("N36").Value = 12.458349
$ws.Range("O36").Value = 0.05436839116396912
$ws.Range("P36").Value = 0.05825596296532109
$ws.Range("Q36").Value = 53.076383147577
$ws.Range("R36").Value = 318.458298885462
$ws.Range("S36").Value = 0.001730085576338225
$ws.Range("T36").Value = 0.001300411289173996

$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 12.780919
$ws.Range("H37").Value = 25.561838
$ws.Range("I37").Value = 0.03182153341857176
$ws.Range("J37").Value = 0.02232237221704002
$ws.Range("K37").Value = 2
$ws.Range("L37").Value = 1
$ws.Range("M37").Value = 11.8826895
$ws.Range("N37").Value = 23.765379
$ws.Range("O37").Value = 0.1555686176754212
$ws.Range("P37").Value = 0.1111282914678999
$ws.Range("Q37").Value = 151.8716920016505
$ws.Range("R37").Value = 607.486768006602
$ws.Range("S37").Value = 0.004950431966239429
$ws.Range("T37").Value = 0.002480647085990175

Write-Host "done"